$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap header row (row 2, C2:G2) with the row-label column (B3:B7) ---
# Before: C2:G2 = "C1".."C5"   and   B3:B7 = "wrt C1".."wrt C5"
# After:  C2:G2 = "wrt C1".."wrt C5"   and   B3:B7 = "C1".."C5"
$ws.Range("C2").Value = "wrt C1"
$ws.Range("D2").Value = "wrt C2"
$ws.Range("E2").Value = "wrt C3"
$ws.Range("F2").Value = "wrt C4"
$ws.Range("G2").Value = "wrt C5"

$ws.Range("B3").Value = "C1"
$ws.Range("B4").Value = "C2"
$ws.Range("B5").Value = "C3"
$ws.Range("B6").Value = "C4"
$ws.Range("B7").Value = "C5"

# --- Reset the pairwise-comparison matrix to all 1s (drop the E3 formula too) ---
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1

$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1

# --- Update the active selection shown in the sheet view ---
$ws.Range("I6").Select() | Out-Null

# --- Remove the inserted picture from the sheet ---
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
  $ws.Shapes.Item($i).Delete() | Out-Null
}
